# Update scripts with new TPM values
# This mirrors a re-run of the NATMI ligand-receptor pipeline with updated
# expression (TPM) inputs, which changes the derived expression/specificity
# columns for rows 2-4 (G,H,I,J,M,N,Q,R,S,T) while leaving identifiers,
# counts/detection-rate flags (E,F,K,L,O,P) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ECs -> Bsg/Sele -> ECs)
$ws.Range("G2").Value = 21.50919633333334
$ws.Range("H2").Value = 64.52758900000001
$ws.Range("I2").Value = 0.1832723264758264
$ws.Range("J2").Value = 0.1832723264758264
$ws.Range("M2").Value = 0.5273163333333333
$ws.Range("N2").Value = 1.581949
$ws.Range("Q2").Value = 11.34215054344011
$ws.Range("R2").Value = 102.079354890961
$ws.Range("S2").Value = 0.1832723264758264
$ws.Range("T2").Value = 0.1832723264758264

# Row 3 (FAPs -> Bsg/Sele -> ECs)
$ws.Range("I3").Value = 0.5927317426910698
$ws.Range("J3").Value = 0.5927317426910698
$ws.Range("M3").Value = 0.5273163333333333
$ws.Range("N3").Value = 1.581949
$ws.Range("Q3").Value = 36.68231198213366
$ws.Range("R3").Value = 330.140807839203
$ws.Range("S3").Value = 0.5927317426910698
$ws.Range("T3").Value = 0.5927317426910698

# Row 4 (MuSCs -> Bsg/Sele -> ECs)
$ws.Range("G4").Value = 26.28859766666667
$ws.Range("H4").Value = 78.86579300000001
$ws.Range("I4").Value = 0.2239959308331038
$ws.Range("J4").Value = 0.2239959308331038
$ws.Range("M4").Value = 0.5273163333333333
$ws.Range("N4").Value = 1.581949
$ws.Range("Q4").Value = 13.86240693006189
$ws.Range("R4").Value = 124.761662370557
$ws.Range("S4").Value = 0.2239959308331038
$ws.Range("T4").Value = 0.2239959308331038
